# Update the "number of views/sales" figures (column F) on the
# "展览" and "全部类型" sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows 2,4,5,7,11) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 305
$ws1.Range("F4").Value = 8086
$ws1.Range("F5").Value = 5894
$ws1.Range("F7").Value = 91
$ws1.Range("F11").Value = 420

# --- Sheet "全部类型" (rows 2,4,5,7,15) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 305
$ws4.Range("F4").Value = 8086
$ws4.Range("F5").Value = 5894
$ws4.Range("F7").Value = 91
$ws4.Range("F15").Value = 420
